$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 2
    "F2" = 0.6666666666666666
    "G2" = 0.2988413333333333
    "H2" = 0.896524
    "I2" = 0.3632971504731247
    "J2" = 0.3632971504731246
    "K2" = 3
    "L2" = 1
    "M2" = 0.1944653333333334
    "N2" = 0.583396
    "O2" = 0.04942840076761122
    "P2" = 0.04942840076761121
    "Q2" = 0.05811427950044445
    "R2" = 0.523028515504
    "S2" = 0.01795719715131677
    "T2" = 0.01795719715131676
    "E3" = 2
    "F3" = 0.6666666666666666
    "G3" = 0.2988413333333333
    "H3" = 0.896524
    "I3" = 0.3632971504731247
    "J3" = 0.3632971504731246
    "O3" = 0.04304638286515546
    "P3" = 0.04304638286515546
    "Q3" = 0.05061077207555556
    "R3" = 0.45549694868
    "S3" = 0.01563862823308612
    "T3" = 0.01563862823308611
    "E4" = 2
    "F4" = 0.6666666666666666
    "G4" = 0.2988413333333333
    "H4" = 0.896524
    "I4" = 0.3632971504731247
    "J4" = 0.3632971504731246
    "O4" = 0.9075252163672334
    "P4" = 0.9075252163672333
    "Q4" = 1.067001425468444
    "R4" = 9.603012829215999
    "S4" = 0.3297013250887218
    "T4" = 0.3297013250887217
    "I5" = 0.4682720202225272
    "J5" = 0.4682720202225272
    "K5" = 3
    "L5" = 1
    "M5" = 0.1944653333333334
    "N5" = 0.583396
    "O5" = 0.04942840076761122
    "P5" = 0.04942840076761121
    "Q5" = 0.07490642585555556
    "R5" = 0.6741578327000001
    "S5" = 0.02314593708381802
    "T5" = 0.02314593708381802
    "I6" = 0.4682720202225272
    "J6" = 0.4682720202225272
    "O6" = 0.04304638286515546
    "P6" = 0.04304638286515546
    "S6" = 0.02015741666753873
    "T6" = 0.02015741666753873
    "I7" = 0.4682720202225272
    "J7" = 0.4682720202225272
    "O7" = 0.9075252163672334
    "P7" = 0.9075252163672333
    "S7" = 0.4249686664711705
    "T7" = 0.4249686664711704
    "I8" = 0.1684308293043481
    "J8" = 0.1684308293043481
    "K8" = 3
    "L8" = 1
    "M8" = 0.1944653333333334
    "N8" = 0.583396
    "O8" = 0.04942840076761122
    "P8" = 0.04942840076761121
    "Q8" = 0.02694278300266667
    "R8" = 0.242485047024
    "S8" = 0.008325266532476436
    "T8" = 0.008325266532476434
    "I9" = 0.1684308293043481
    "J9" = 0.1684308293043481
    "O9" = 0.04304638286515546
    "P9" = 0.04304638286515546
    "Q9" = 0.02346402745333334
    "S9" = 0.007250337964530616
    "T9" = 0.007250337964530615
    "I10" = 0.1684308293043481
    "J10" = 0.1684308293043481
    "O10" = 0.9075252163672334
    "P10" = 0.9075252163672333
    "Q10" = 0.4946802768106666
    "R10" = 4.452122491296
    "S10" = 0.1528552248073411
    "T10" = 0.1528552248073411
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
